$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout had 5 rows (with two blank spacer rows). The new layout
# is a clean 3-row x 5-col table, so drop the now-unused rows first.
$ws.Rows("4:5").Delete()

# Make sure numeric-looking labels ("2019", "1,177,951", ...) are stored
# as text, matching the shared-string table in the target workbook.
$ws.Range("A1:E3").NumberFormat = "@"

# Fill column by column (top-to-bottom, then next column) so the shared
# string table is built up in the same order as the source data extract.
$ws.Range("A1").Value = "Year"
$ws.Range("A2").Value = "2019"
$ws.Range("A3").Value = "2018"

$ws.Range("B1").Value = "Company cars"
$ws.Range("B2").Value = "1,177,951"
$ws.Range("B3").Value = "1,162,063"

$ws.Range("C1").Value = "Total by rental car/private car"
$ws.Range("C2").Value = "537,721"
$ws.Range("C3").Value = "590,929"

$ws.Range("D1").Value = "Total by rail"
$ws.Range("D2").Value = "45,722"
$ws.Range("D3").Value = "47,804"

$ws.Range("E1").Value = "Total by air*"
$ws.Range("E2").Value = "2,103,706"
$ws.Range("E3").Value = "867,678"
